# Insert two new italic "BodyText" commentary notes right after the
# paragraph ending in "...interannual time scales." and before the
# existing "(one paragraph discussing models)" paragraph.

$d = $word.ActiveDocument

# Locate the end of the sentence that the new paragraphs must follow.
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "However, we now have clear evidence that it is not realistic to expect that a constant allocation of photosynthate will be allocated to woody growth on either intraanuual or interannual time scales.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor sentence for commentary insertion."
}

# Collapse to the end of the matched sentence (right before the paragraph mark).
$anchor.Collapse(0)

# --- First new paragraph: "(active vs passive allocation)" ---
$anchor.InsertParagraphAfter() | Out-Null
$anchor.Move(1, 1) | Out-Null
$anchor.InsertAfter("(active vs passive allocation)")

# --- Second new paragraph: "(implications for tree-ring studies)" ---
$anchor.Collapse(0)
$anchor.InsertParagraphAfter() | Out-Null
$anchor.Move(1, 1) | Out-Null
$anchor.InsertAfter("(implications for tree-ring studies)")

Write-Output "Inserted commentary paragraphs."
